$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5146761536598206
$ws.Range("B1").Value = 0.6105854511260986
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 1.470643639564514
$ws.Range("E1").Value = 0.8930754065513611
